$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 7045
$ws.Range("I62").Value = 6039.2856
$ws.Range("K62").Value = 6039.2856
$ws.Range("M62").Value = -5415.2856
$ws.Range("H65").Value = 7045
$ws.Range("I65").Value = 6039.2856
$ws.Range("K65").Value = 30196.428
$ws.Range("M65").Value = -27076.428
$ws.Range("H70").Value = 102838.86
$ws.Range("I70").Value = 3250
$ws.Range("K70").Value = 9750
$ws.Range("M70").Value = -9480
$ws.Range("H73").Value = 102838.86
$ws.Range("I73").Value = 3250
$ws.Range("K73").Value = 9750
$ws.Range("M73").Value = -8814
$ws.Range("H96").Value = 1266.2106
$ws.Range("J96").Value = 1203.3334
$ws.Range("L96").Value = 3610.0002
$ws.Range("N96").Value = -6356.0002
$ws.Range("H112").Value = 2453.3
$ws.Range("J112").Value = 2731.1428
$ws.Range("L112").Value = 8193.428400000001
$ws.Range("N112").Value = -10409.4284
$ws.Range("H125").Value = 1381
$ws.Range("I125").Value = 1174
$ws.Range("J125").Value = 1795
$ws.Range("K125").Value = 10566
$ws.Range("L125").Value = 16155
$ws.Range("M125").Value = -8106
$ws.Range("N125").Value = -21075
$ws.Range("H137").Value = 13840.917
$ws.Range("I137").Value = 15954.777
$ws.Range("K137").Value = 47864.331
$ws.Range("M137").Value = -45314.331
$ws.Range("H138").Value = 11356.191
$ws.Range("I138").Value = 10067.3
$ws.Range("J138").Value = 11704.541
$ws.Range("K138").Value = 30201.9
$ws.Range("L138").Value = 35113.623
$ws.Range("M138").Value = -25061.9
$ws.Range("N138").Value = -45393.623

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 24595.516
$ws.Range("I32").Value = 21223.65
$ws.Range("K32").Value = 21223.65
$ws.Range("M32").Value = -20936.65
$ws.Range("H45").Value = 2696.2
$ws.Range("I45").Value = 2696.2
$ws.Range("K45").Value = 2696.2
$ws.Range("M45").Value = -2319.2
$ws.Range("H61").Value = 1638
$ws.Range("I61").Value = 1512.6666
$ws.Range("J61").Value = 2014
$ws.Range("K61").Value = 1512.6666
$ws.Range("L61").Value = 2014
$ws.Range("M61").Value = -1300.6666
$ws.Range("N61").Value = -2438
$ws.Range("H74").Value = 3024.4707
$ws.Range("I74").Value = 1090.4
$ws.Range("K74").Value = 1090.4
$ws.Range("M74").Value = -216.4000000000001
$ws.Range("H77").Value = 3024.4707
$ws.Range("I77").Value = 1090.4
$ws.Range("K77").Value = 5452
$ws.Range("M77").Value = -1084
$ws.Range("H132").Value = 2820.6667
$ws.Range("I132").Value = 2346.4
$ws.Range("K132").Value = 7039.200000000001
$ws.Range("M132").Value = -4509.200000000001
$ws.Range("H136").Value = 1638
$ws.Range("I136").Value = 1512.6666
$ws.Range("J136").Value = 2014
$ws.Range("K136").Value = 4537.9998
$ws.Range("L136").Value = 6042
$ws.Range("M136").Value = -1987.9998
$ws.Range("N136").Value = -11142

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 15000
$ws.Range("J9").Value = 15000
$ws.Range("L9").Value = 15000
$ws.Range("N9").Value = -15336
$ws.Range("H26").Value = 21996.25
$ws.Range("I26").Value = 21996.25
$ws.Range("K26").Value = 21996.25
$ws.Range("M26").Value = -21704.25
$ws.Range("H44").Value = 39999
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 39999
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 39999
$ws.Range("M44").ClearContents()
$ws.Range("N44").Value = -40993
$ws.Range("H86").Value = 6678.9165
$ws.Range("I86").Value = 6642.5713
$ws.Range("J86").Value = 6729.8
$ws.Range("K86").Value = 6642.5713
$ws.Range("L86").Value = 6729.8
$ws.Range("M86").Value = -5519.5713
$ws.Range("N86").Value = -8975.799999999999
$ws.Range("H89").Value = 6678.9165
$ws.Range("I89").Value = 6642.5713
$ws.Range("J89").Value = 6729.8
$ws.Range("K89").Value = 33212.85649999999
$ws.Range("L89").Value = 33649
$ws.Range("M89").Value = -27596.85649999999
$ws.Range("N89").Value = -44881

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4180.615
$ws.Range("I31").Value = 2555.8125
$ws.Range("K31").Value = 2555.8125
$ws.Range("M31").Value = -2260.8125
$ws.Range("H34").Value = 4180.615
$ws.Range("I34").Value = 2555.8125
$ws.Range("K34").Value = 2555.8125
$ws.Range("M34").Value = -2353.8125
$ws.Range("H107").Value = 722.7895
$ws.Range("I107").Value = 403.83334
$ws.Range("J107").Value = 1269.5714
$ws.Range("K107").Value = 403.83334
$ws.Range("L107").Value = 1269.5714
$ws.Range("M107").Value = 1516.16666
$ws.Range("N107").Value = -5109.5714
$ws.Range("H132").Value = 2307.3333
$ws.Range("I132").Value = 2207.4138
$ws.Range("J132").Value = 2721.2856
$ws.Range("K132").Value = 6622.241399999999
$ws.Range("L132").Value = 8163.8568
$ws.Range("M132").Value = -4092.241399999999
$ws.Range("N132").Value = -13223.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 702.8333
$ws.Range("J114").Value = 729.25
$ws.Range("L114").Value = 2187.75
$ws.Range("N114").Value = -8695.75
$ws.Range("H131").Value = 3756
$ws.Range("I131").Value = 1420
$ws.Range("K131").Value = 4260
$ws.Range("M131").Value = 780
$ws.Range("H139").Value = 18066.334
$ws.Range("I139").Value = 29200
$ws.Range("K139").Value = 87600
$ws.Range("M139").Value = -82460
$ws.Range("H140").Value = 4867.2
$ws.Range("I140").Value = 4334
$ws.Range("J140").Value = 7000
$ws.Range("K140").Value = 13002
$ws.Range("L140").Value = 21000
$ws.Range("M140").Value = -7822
$ws.Range("N140").Value = -31360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 9599.200000000001
$ws.Range("I80").Value = 8999
$ws.Range("K80").Value = 8999
$ws.Range("M80").Value = -8001
$ws.Range("H83").Value = 9599.200000000001
$ws.Range("I83").Value = 8999
$ws.Range("K83").Value = 44995
$ws.Range("M83").Value = -40003

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 8708.833000000001
$ws.Range("I7").Value = 8450.6
$ws.Range("K7").Value = 8450.6
$ws.Range("M7").Value = -8338.6
$ws.Range("H22").Value = 3500
$ws.Range("I22").Value = 2000
$ws.Range("J22").Value = 5000
$ws.Range("K22").Value = 2000
$ws.Range("L22").Value = 5000
$ws.Range("M22").Value = -1705
$ws.Range("N22").Value = -5590
$ws.Range("H27").Value = 3500
$ws.Range("I27").Value = 2000
$ws.Range("J27").Value = 5000
$ws.Range("K27").Value = 2000
$ws.Range("L27").Value = 5000
$ws.Range("M27").Value = -1893
$ws.Range("N27").Value = -5214
$ws.Range("H46").Value = 3064.9473
$ws.Range("J46").Value = 3171.8462
$ws.Range("L46").Value = 3171.8462
$ws.Range("N46").Value = -3547.8462
$ws.Range("H100").Value = 4225.375
$ws.Range("I100").Value = 3467
$ws.Range("J100").Value = 6500.5
$ws.Range("K100").Value = 3467
$ws.Range("L100").Value = 6500.5
$ws.Range("M100").Value = -2926
$ws.Range("N100").Value = -7582.5
$ws.Range("H126").Value = 8708.833000000001
$ws.Range("I126").Value = 8450.6
$ws.Range("K126").Value = 25351.8
$ws.Range("M126").Value = -22881.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6552.2
$ws.Range("I62").Value = 5665.6665
$ws.Range("J62").Value = 6773.8335
$ws.Range("K62").Value = 5665.6665
$ws.Range("L62").Value = 6773.8335
$ws.Range("M62").Value = -5041.6665
$ws.Range("N62").Value = -8021.8335
$ws.Range("H65").Value = 6552.2
$ws.Range("I65").Value = 5665.6665
$ws.Range("J65").Value = 6773.8335
$ws.Range("K65").Value = 28328.3325
$ws.Range("L65").Value = 33869.1675
$ws.Range("M65").Value = -25208.3325
$ws.Range("N65").Value = -40109.1675
$ws.Range("H81").Value = 4277.643
$ws.Range("J81").Value = 5999
$ws.Range("L81").Value = 11998
$ws.Range("N81").Value = -14120
$ws.Range("H84").Value = 4277.643
$ws.Range("J84").Value = 5999
$ws.Range("L84").Value = 59990
$ws.Range("N84").Value = -70598
$ws.Range("H96").Value = 2771.3333
$ws.Range("I96").Value = 2771.3333
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 2771.3333
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -1398.3333
$ws.Range("N96").ClearContents()
$ws.Range("H100").Value = 2579
$ws.Range("I100").Value = 2579
$ws.Range("K100").Value = 5158
$ws.Range("M100").Value = -4617
$ws.Range("H132").Value = 2022.3448
$ws.Range("I132").Value = 2136.8696
$ws.Range("J132").Value = 1583.3334
$ws.Range("K132").Value = 6410.6088
$ws.Range("L132").Value = 4750.0002
$ws.Range("M132").Value = -3880.6088
$ws.Range("N132").Value = -9810.0002
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
